$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.136.06"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -2.80%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.868.55"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -2.17%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'0.9986"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'306.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -1.83%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'0.9987"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -0.01%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.5121"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +2.60%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.3747"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  -1.55%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  -2.07%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.8871"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -2.54%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'20.62"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  -3.24%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'1.879.24"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +0.68%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.07521"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -1.55%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'5.311"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -3.06%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'89.17"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -3.67%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.9986"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -0.02%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'0.000008471"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -2.88%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'14.09"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -3.78%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'0.9996"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +0.13%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'27.166.59"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -2.82%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'5.052"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -2.21%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'2.112.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  -0.58%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'10.55"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -2.94%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'6.483"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -1.78%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'1.850"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  +0.30%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'149.48"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -2.25%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'17.94"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -2.51%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'2.104"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -5.22%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'112.82"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -1.90%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'4.737"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -3.43%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'4.684"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -2.63%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'0.09016"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +0.29%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'0.05126"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -2.89%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'3.085"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -3.61%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'1.162"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -5.68%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'0.7350"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -5.70%  "
$ws.Range("E36").ClearFormats()
$ws.Range("E37").Value = "'  -1.76%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'2.505"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -4.47%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'3.052"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -0.24%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'1.079"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -1.10%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.5338"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -3.93%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'6.579"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -3.39%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'116.61"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +2.29%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'8.330"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -2.11%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.1473"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  -2.82%  "
$ws.Range("E45").ClearFormats()
$ws.Range("B46").Value = "'Decentraland"
$ws.Range("B46").ClearFormats()
$ws.Range("C46").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("C46").ClearFormats()
$ws.Range("D46").Value = "'0.4630"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -4.05%  "
$ws.Range("E46").ClearFormats()
$ws.Range("B47").Value = "'PaxDollar"
$ws.Range("B47").ClearFormats()
$ws.Range("C47").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("C47").ClearFormats()
$ws.Range("D47").Value = "'0.9983"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -0.03%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'10.05"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  -5.53%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'1.569"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -4.39%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'64.40"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -4.30%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'36.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -1.57%  "
$ws.Range("E51").ClearFormats()
